$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.086.60'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').Value = '  -2.85%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.711.39'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').Value = '  -3.27%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.96'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').Value = '  -6.23%  '

$ws.Range('E6').Value = '  +0.11%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4725'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').Value = '  +5.17%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3417'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').Value = '  -3.90%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '41.92'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').Value = '  -0.28%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07251'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').Value = '  -2.38%  '

$ws.Range('E11').Value = '  -5.82%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').Value = '  +0.02%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.78'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').Value = '  -5.70%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.851'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').Value = '  -2.95%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.713.60'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').Value = '  -3.34%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.853'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').Value = '  -5.45%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.07'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').Value = '  -4.45%  '

$ws.Range('E18').Value = '  -2.20%  '

$ws.Range('E19').Value = '  -1.35%  '

$ws.Range('E20').Value = '  +0.07%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.49'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').Value = '  -3.70%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.605'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').Value = '  -3.14%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.121.84'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').Value = '  -2.89%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.86'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').Value = '  -3.69%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.106'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').Value = '  +0.02%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.78'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').Value = '  -4.18%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.44'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').Value = '  -4.62%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.910.73'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').Value = '  -3.31%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.079'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').Value = '  -3.87%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.99'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').Value = '  -4.57%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.010'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').Value = '  -8.87%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09143'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').Value = '  -0.48%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.587'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').Value = '  -1.99%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.300'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').Value = '  -5.86%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02200'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').Value = '  -3.99%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05806'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').Value = '  -4.88%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.93'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').Value = '  -7.90%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.1987'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').Value = '  -5.47%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.729'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').Value = '  -4.72%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.399'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').Value = '  +0.36%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5860'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').Value = '  -7.40%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.110'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').Value = '  -6.05%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.443'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').Value = '  -5.76%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.46'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').Value = '  -6.07%  '

$ws.Range('E45').Value = '  -5.04%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5634'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').Value = '  -4.37%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '116.84'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').Value = '  -4.79%  '

$ws.Range('E48').Value = '  -6.41%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06639'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').Value = '  -3.83%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.080'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').Value = '  -5.15%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').Value = '  +0.18%  '
